# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates the "K" column (spreadsheet column G) values for the
# save-data rows, reflecting a recalculation of the K statistic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 2
    4  = 2
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    12 = 4
    13 = 2
    14 = 3
    15 = 1
    16 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
